# Update attendance/view counts (column F) on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets to the newly-generated numbers.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 -----------------------------------------------------
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Cells.Item(2, 6).Value = 53
$wsExhibit.Cells.Item(4, 6).Value = 1586
$wsExhibit.Cells.Item(5, 6).Value = 287
$wsExhibit.Cells.Item(6, 6).Value = 78
$wsExhibit.Cells.Item(7, 6).Value = 1819
$wsExhibit.Cells.Item(8, 6).Value = 10256
$wsExhibit.Cells.Item(9, 6).Value = 176
$wsExhibit.Cells.Item(11, 6).Value = 260
$wsExhibit.Cells.Item(14, 6).Value = 7112
$wsExhibit.Cells.Item(15, 6).Value = 1105
$wsExhibit.Cells.Item(16, 6).Value = 668
$wsExhibit.Cells.Item(17, 6).Value = 62
$wsExhibit.Cells.Item(19, 6).Value = 242

# --- Sheet 4: 全部类型 --------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Cells.Item(2, 6).Value = 53
$wsAll.Cells.Item(4, 6).Value = 1586
$wsAll.Cells.Item(5, 6).Value = 287
$wsAll.Cells.Item(7, 6).Value = 78
$wsAll.Cells.Item(8, 6).Value = 1819
$wsAll.Cells.Item(11, 6).Value = 10256
$wsAll.Cells.Item(12, 6).Value = 176
$wsAll.Cells.Item(14, 6).Value = 260
$wsAll.Cells.Item(17, 6).Value = 7112
$wsAll.Cells.Item(18, 6).Value = 1105
$wsAll.Cells.Item(19, 6).Value = 668
$wsAll.Cells.Item(20, 6).Value = 62
$wsAll.Cells.Item(22, 6).Value = 242
